$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace record_id placeholder values (a,b,c,d,e) with the real error message
# now produced for records that do not exist in the study.
$ws.Range("A2").Value = "Error: record does not exist in study"
$ws.Range("A3").Value = "Error: record does not exist in study"
$ws.Range("A4").Value = "Error: record does not exist in study"
$ws.Range("A5").Value = "Error: record does not exist in study"
$ws.Range("A6").Value = "Error: record does not exist in study"

# Widen column A to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 32.17

# Reset the view: drop the frozen/scrolled-to "E1" top-left cell and move
# the active selection down to B11 (a single cell instead of the old I2:J6 block).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()
